# Apply the edits described by the commit: "Updated `kundur_reg.xlsx`; minor
# tweaks to `s_update_var`".
#
# Real semantic changes (everything else in the diff is float round-trip /
# Excel-version metadata noise produced by re-saving the workbook):
#
#   1. REPCA1 sheet: insert a new parameter column "Kc" right before the
#      existing "emax" column (i.e. before column T), filled with 1 for the
#      one data row, and tweak a handful of other REPCA1 values.
#   2. Toggler sheet: row for idx=2 now has u=1 (was 0).
#   3. REPCA1 becomes the active sheet/tab, with the view scrolled so column S
#      is the left-most visible column and the active cell at AL3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. REPCA1 ("REPCA1" sheet) -- insert new "Kc" column and update values
# ---------------------------------------------------------------------------
$repca = $wb.Worksheets.Item("REPCA1")

# Insert a new column before column T (20th column), shifting emax..Dup right.
$repca.Columns.Item(20).Insert()

# New column header + value.
$repca.Cells.Item(1, 20).Value = "Kc"
$repca.Cells.Item(2, 20).Value = 1

# Other REPCA1 value tweaks (columns after the insert point are shifted by one
# relative to the original file).
$repca.Cells.Item(2, 9).Value = 1        # I2  VCFlag:  0    -> 1
$repca.Cells.Item(2, 10).Value = 1       # J2  RefFlag: 0    -> 1
$repca.Cells.Item(2, 17).Value = 0.98    # Q2  Vfrz:    0.8  -> 0.98
$repca.Cells.Item(2, 23).Value = -0.02   # W2  dbd1:    -0.1 -> -0.02
$repca.Cells.Item(2, 24).Value = 0.02    # X2  dbd2:    0.1  -> 0.02
$repca.Cells.Item(2, 35).Value = -999    # AI2 Pmin:    0    -> -999
$repca.Cells.Item(2, 37).Value = 10      # AK2 Ddn:     0.05 -> 10
$repca.Cells.Item(2, 38).Value = 10      # AL2 Dup:     0.05 -> 10

# ---------------------------------------------------------------------------
# 2. Toggler sheet -- enable (u=1) the second toggler row (idx=2)
# ---------------------------------------------------------------------------
$toggler = $wb.Worksheets.Item("Toggler")
$toggler.Cells.Item(3, 3).Value = 1      # C3  u: 0 -> 1

# ---------------------------------------------------------------------------
# 3. Make REPCA1 the active sheet with the saved scroll/selection state
# ---------------------------------------------------------------------------
$repca.Activate()
$repca.Application.ActiveWindow.ScrollColumn = 19  # topLeftCell="S1"
$repca.Range("AL3").Select()
